# chore: update Sheets via scheduled runner
# Refresh cached market-price / profit figures (columns H-N) across the
# per-job leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I32").Value = 849.5
$ws.Range("J32").Value = 959.38464
$ws.Range("K32").Value = 849.5
$ws.Range("L32").Value = 959.38464
$ws.Range("M32").Value = -523.5
$ws.Range("N32").Value = -1611.38464

$ws.Range("H39").Value = 458.8125
$ws.Range("I39").Value = 60.583332
$ws.Range("J39").Value = 1653.5
$ws.Range("K39").Value = 181.749996
$ws.Range("L39").Value = 4960.5
$ws.Range("M39").Value = 114.250004
$ws.Range("N39").Value = -5552.5

$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 3500
$ws.Range("L51").Value = 3500
$ws.Range("N51").Value = -4468

$ws.Range("H74").Value = 3286.889
$ws.Range("I74").Value = 3136.182
$ws.Range("J74").Value = 3523.7144
$ws.Range("K74").Value = 3136.182
$ws.Range("L74").Value = 3523.7144
$ws.Range("M74").Value = -2200.182
$ws.Range("N74").Value = -5395.7144

$ws.Range("H77").Value = 3286.889
$ws.Range("I77").Value = 3136.182
$ws.Range("J77").Value = 3523.7144
$ws.Range("K77").Value = 15680.91
$ws.Range("L77").Value = 17618.572
$ws.Range("M77").Value = -11000.91
$ws.Range("N77").Value = -26978.572

$ws.Range("H103").Value = 4007167.2
$ws.Range("I103").Value = 8586044
$ws.Range("J103").Value = 650
$ws.Range("K103").Value = 25758132
$ws.Range("L103").Value = 1950
$ws.Range("M103").Value = -25757546
$ws.Range("N103").Value = -3122

$ws.Range("H131").Value = 2526.3635
$ws.Range("I131").Value = 758
$ws.Range("J131").Value = 4000
$ws.Range("K131").Value = 2274
$ws.Range("L131").Value = 12000
$ws.Range("M131").Value = 2766
$ws.Range("N131").Value = -22080

$ws.Range("H137").Value = 6454845
$ws.Range("I137").Value = 2179.6
$ws.Range("J137").Value = 12504219
$ws.Range("K137").Value = 6538.799999999999
$ws.Range("L137").Value = 37512657
$ws.Range("M137").Value = -3988.799999999999
$ws.Range("N137").Value = -37517757

$ws.Range("H138").Value = 7579052
$ws.Range("I138").Value = 1296.88
$ws.Range("J138").Value = 31259538
$ws.Range("K138").Value = 3890.64
$ws.Range("L138").Value = 93778614
$ws.Range("M138").Value = 1249.36
$ws.Range("N138").Value = -93788894

$ws.Range("H141").Value = 1574.5
$ws.Range("I141").Value = 1435.3636
$ws.Range("J141").Value = 3105
$ws.Range("K141").Value = 4306.0908
$ws.Range("L141").Value = 9315
$ws.Range("M141").Value = 873.9092000000001
$ws.Range("N141").Value = -19675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14149.289
$ws.Range("I32").Value = 16609.75
$ws.Range("K32").Value = 16609.75
$ws.Range("M32").Value = -16322.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 532.9091
$ws.Range("I20").Value = 530.25
$ws.Range("J20").Value = 540
$ws.Range("K20").Value = 530.25
$ws.Range("L20").Value = 540
$ws.Range("M20").Value = -283.25
$ws.Range("N20").Value = -1034

$ws.Range("H105").Value = 3068.75
$ws.Range("I105").Value = 1855.5555
$ws.Range("K105").Value = 1855.5555
$ws.Range("M105").Value = -108.5554999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9263944
$ws.Range("I31").Value = 10097.214
$ws.Range("J31").Value = 15152756
$ws.Range("K31").Value = 10097.214
$ws.Range("L31").Value = 15152756
$ws.Range("M31").Value = -9802.214
$ws.Range("N31").Value = -15153346

$ws.Range("H34").Value = 9263944
$ws.Range("I34").Value = 10097.214
$ws.Range("J34").Value = 15152756
$ws.Range("K34").Value = 10097.214
$ws.Range("L34").Value = 15152756
$ws.Range("M34").Value = -9895.214
$ws.Range("N34").Value = -15153160

$ws.Range("H99").Value = 1471.4286
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 1460
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 1460
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -4456

$ws.Range("H122").Value = 1293.88
$ws.Range("I122").Value = 1295.2778
$ws.Range("J122").Value = 1290.2858
$ws.Range("K122").Value = 3885.8334
$ws.Range("L122").Value = 3870.8574
$ws.Range("M122").Value = -1435.8334
$ws.Range("N122").Value = -8770.857400000001

$ws.Range("H126").Value = 1471.4286
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 1460
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 4380
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -9320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1389085.8
$ws.Range("I2").Value = 203.33333
$ws.Range("J2").Value = 3472409.5
$ws.Range("K2").Value = 1219.99998
$ws.Range("L2").Value = 20834457
$ws.Range("M2").Value = -1106.99998
$ws.Range("N2").Value = -20834683

$ws.Range("H5").Value = 514.29034
$ws.Range("I5").Value = 368.875
$ws.Range("J5").Value = 1012.8571
$ws.Range("K5").Value = 1106.625
$ws.Range("L5").Value = 3038.5713
$ws.Range("M5").Value = -994.625
$ws.Range("N5").Value = -3262.5713

$ws.Range("H17").Value = 712.75
$ws.Range("I17").Value = 475.5
$ws.Range("J17").Value = 950
$ws.Range("K17").Value = 1426.5
$ws.Range("L17").Value = 2850
$ws.Range("M17").Value = -1257.5
$ws.Range("N17").Value = -3188

$ws.Range("H55").Value = 549.9375
$ws.Range("I55").Value = 404
$ws.Range("J55").Value = 559.6667
$ws.Range("K55").Value = 1212
$ws.Range("L55").Value = 1679.0001
$ws.Range("M55").Value = -1035
$ws.Range("N55").Value = -2033.0001

$ws.Range("H135").Value = 514.29034
$ws.Range("I135").Value = 368.875
$ws.Range("J135").Value = 1012.8571
$ws.Range("K135").Value = 3319.875
$ws.Range("L135").Value = 9115.713899999999
$ws.Range("M135").Value = -784.875
$ws.Range("N135").Value = -14185.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 3834.8333
$ws.Range("I22").Value = 2600
$ws.Range("J22").Value = 10009
$ws.Range("K22").Value = 2600
$ws.Range("L22").Value = 10009
$ws.Range("M22").Value = -2071
$ws.Range("N22").Value = -11067

$ws.Range("H80").Value = 11907538
$ws.Range("I80").Value = 30305658
$ws.Range("J80").Value = 2871.4119
$ws.Range("K80").Value = 30305658
$ws.Range("L80").Value = 2871.4119
$ws.Range("M80").Value = -30304660
$ws.Range("N80").Value = -4867.4119

$ws.Range("H83").Value = 11907538
$ws.Range("I83").Value = 30305658
$ws.Range("J83").Value = 2871.4119
$ws.Range("K83").Value = 151528290
$ws.Range("L83").Value = 14357.0595
$ws.Range("M83").Value = -151523298
$ws.Range("N83").Value = -24341.0595

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2368.4443
$ws.Range("I16").Value = 2368.4443
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2368.4443
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2198.4443
$ws.Range("N16").ClearContents()

$ws.Range("H82").Value = 2089.818
$ws.Range("I82").Value = 1501.6
$ws.Range("J82").Value = 2580
$ws.Range("K82").Value = 1501.6
$ws.Range("L82").Value = 2580
$ws.Range("M82").Value = -1140.6
$ws.Range("N82").Value = -3302

$ws.Range("H85").Value = 2089.818
$ws.Range("I85").Value = 1501.6
$ws.Range("J85").Value = 2580
$ws.Range("K85").Value = 1501.6
$ws.Range("L85").Value = 2580
$ws.Range("M85").Value = -253.5999999999999
$ws.Range("N85").Value = -5076

$ws.Range("H100").Value = 1981.6
$ws.Range("I100").Value = 1500
$ws.Range("K100").Value = 1500
$ws.Range("M100").Value = -959

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7585.4443
$ws.Range("J74").Value = 7962.5
$ws.Range("L74").Value = 7962.5
$ws.Range("N74").Value = -9834.5

$ws.Range("H77").Value = 7585.4443
$ws.Range("J77").Value = 7962.5
$ws.Range("L77").Value = 23887.5
$ws.Range("N77").Value = -33247.5

$ws.Range("H100").Value = 5343
$ws.Range("I100").Value = 7245
$ws.Range("J100").Value = 2490
$ws.Range("K100").Value = 14490
$ws.Range("L100").Value = 4980
$ws.Range("M100").Value = -13949
$ws.Range("N100").Value = -6062
